# issue #5: stock data output to json file
# Adds a "property_category" column (with constant value "stock") to the
# 股票 (stock) sheet, between the existing "total" and "date" columns, and
# fixes a stray apostrophe in one of the stock name values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the old column H ("date"); this pushes the
# existing H/I/J (date / legislator_name / legislator_id) columns one to
# the right, copying their formatting along the way.
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Fill every data row (rows 2-15) with the constant category value "stock".
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Fix the stray apostrophe in the "勤益" stock name (row 12, column B).
$ws.Cells.Item(12, 2).Value = "勤益"
